$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F7").Value = 7610
$ws.Range("F9").Value = 176
$ws.Range("F10").Value = 6704
$ws.Range("F11").Value = 138
$ws.Range("F12").Value = 288
$ws.Range("F13").Value = 4728
$ws.Range("F17").Value = 4994
$ws.Range("F20").Value = 290
$ws.Range("F26").Value = 8616
$ws.Range("F27").Value = 64
$ws.Range("F28").Value = 35
$ws.Range("F30").Value = 761
$ws.Range("F37").Value = 1770
$ws.Range("F39").Value = 1062
$ws.Range("F41").Value = 4519
$ws.Range("F45").Value = 59
$ws.Range("F46").Value = 11
$ws.Range("F47").Value = 875
$ws.Range("F48").Value = 1184

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F17").Value = 878

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 244

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 244
$ws.Range("F9").Value = 7610
$ws.Range("F11").Value = 176
$ws.Range("F12").Value = 6704
$ws.Range("F13").Value = 138
$ws.Range("F14").Value = 288
$ws.Range("F15").Value = 4729
$ws.Range("F19").Value = 4994
$ws.Range("F22").Value = 290
$ws.Range("F29").Value = 8616
$ws.Range("F30").Value = 64
$ws.Range("F31").Value = 35
$ws.Range("F33").Value = 761
$ws.Range("F38").Value = 1770
$ws.Range("F40").Value = 1062
$ws.Range("F42").Value = 4519
$ws.Range("F46").Value = 59
$ws.Range("F47").Value = 875
$ws.Range("F48").Value = 1184
